$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the typo / update text of existing comment (K32) ---
# This updates the shared string in place (index 85) exactly like the diff.
$ws.Range("K32").Value = "doesn't match +  (typo book)"

# --- New data rows for Chapter 17 (rows 33-40) ---
$rowsData = @(
    @{ Row=33; A=17; B=930;  C="MVAMDScity1";     D=930;  E="Ready"; F=462;  G="X"; H="-"; I="Ready"; J="09.09.2016"; K="esthetic differences"; KHighlight=$false },
    @{ Row=34; A=17; B=931;  C="MVAMDScity2";     D=931;  E="Ready"; F=431;  G="X"; H="-"; I="Ready"; J="09.09.2016"; K="esthetic differences"; KHighlight=$false },
    @{ Row=35; A=17; B=1214; C="MVAmdscarm";      D=1214; E="Ready"; F=1687; G="X"; H="-"; I="Ready"; J="09.09.2016"; K="graphs don't match";  KHighlight=$true  },
    @{ Row=36; A=17; B=933;  C="MVAMDSpooladj";   D=933;  E="Ready"; F=318;  G="X"; H="-"; I="Ready"; J="09.09.2016"; K="esthetic differences"; KHighlight=$false },
    @{ Row=37; A=17; B=932;  C="MVAMDSnonmstart"; D=932;  E="Ready"; F=659;  G="X"; H="-"; I="Ready"; J="09.09.2016"; K="esthetic differences"; KHighlight=$false },
    @{ Row=38; A=17; B=934;  C="MVAnmdscar1";     D=934;  E="Ready"; F=482;  G="X"; H="-"; I="Ready"; J="09.09.2016"; K="ok";                   KHighlight=$false },
    @{ Row=39; A=17; B=935;  C="MVAnmdscar2";     D=935;  E="Ready"; F=660;  G="X"; H="-"; I="Ready"; J="09.09.2016"; K="esthetic differences"; KHighlight=$false },
    @{ Row=40; A=17; B=936;  C="MVAnmdscar3";     D=936;  E="Ready"; F=661;  G="X"; H="-"; I="Ready"; J="09.09.2016"; K="esthetic differences"; KHighlight=$false }
)

foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $kCell = $ws.Cells.Item($row, 11)
    $kCell.Value = $r.K
    if ($r.KHighlight) {
        $kCell.Interior.Color = 0x6600FF
        $kCell.HorizontalAlignment = -4108
        $kCell.VerticalAlignment = -4108
    }
}

# --- New Chapter 18 marker row (row 41, only column A populated) ---
$ws.Cells.Item(41, 1).Value = 18

# --- Update the view: scroll so row 18 is the first unfrozen visible row, and
#     move the active selection to the first empty row after the new data. ---
$ws.Range("A42").Select()
$excel.ActiveWindow.ScrollRow = 18
